$wb = $excel.ActiveWorkbook

# --- GLOBAL RESULTS sheet ---
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")

$ws.Range("C2").Value  = 5.329030371724636
$ws.Range("C3").Value  = 12.453995006296132
$ws.Range("C4").Value  = 0.704518060315964
$ws.Range("C6").Value  = 5.090243222753136
$ws.Range("C7").Value  = 11.904207711875411
$ws.Range("C8").Value  = 0.7677698827414503
$ws.Range("C10").Value = 5.090243222753136
$ws.Range("C11").Value = 11.904207711875411
$ws.Range("C12").Value = 0.7677698827414503
$ws.Range("C14").Value = 5.116284719535766
$ws.Range("C15").Value = 11.964166064410108
$ws.Range("C16").Value = 0.48122601702747003
$ws.Range("C18").Value = 4.992295568767794
$ws.Range("C19").Value = 11.67869148995392
$ws.Range("C20").Value = 0.6933403021326805

# --- LANDING GEARS sheet ---
$ws2 = $wb.Worksheets.Item("LANDING GEARS")
$ws2.Range("C2").Value = 12.30269623209934
